$wb = $excel.ActiveWorkbook

# --- Overview sheet: row for 17232927-...md switches from "Ready for
#     handoff" to "Handed back: in sync with en-US" for both locale columns.
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("E3").Value = "Handed back: in sync with en-US"
$wsOverview.Range("F3").Value = "Handed back: in sync with en-US"

# --- zh-cn sheet: same row's Status flips, the Latest Handback DateTime is
#     refreshed, and the stale-handback Error Detail is cleared now that the
#     handback report was (re)generated.
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("C3").Value = "Handed back: in sync with en-US"
$wsZhCn.Range("K3").Value = "2016-10-27 02:17:26"
$wsZhCn.Range("P3").Value = ""

# --- de-de sheet: mirror of the zh-cn update, with its own handback time.
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("C3").Value = "Handed back: in sync with en-US"
$wsDeDe.Range("K3").Value = "2016-10-27 02:17:43"
$wsDeDe.Range("P3").Value = ""

# --- Error Detail column narrows now that it's no longer holding long text.
$wsZhCn.Columns.Item(16).ColumnWidth = 13.7470531463623
$wsDeDe.Columns.Item(16).ColumnWidth = 13.7470531463623
